{"js": "// 1) Append \" (Kijin)\" to the paragraph that reads \"1:15  PTM\" (Day 2 PTM tutorial line).\nconst body = context.document.body;\nconst ptmResults = body.search(\"PTM\", { matchCase: true, matchWholeWord: false });\nptmResults.load(\"items\");\nawait context.sync();\n\nlet ptmRange = null;\nfor (let i = 0; i < ptmResults.items.length; i++) {\n  const par = ptmResults.items[i].paragraphs.getFirst();\n  par.load(\"text\");\n  await context.sync();\n  // The target line is exactly \"1:15  PTM\" (the other \"PTM\" hit is \"...up to PTM)\").\n  if (/^\\s*1:15\\s+PTM\\s*$/.test(par.text)) {\n    ptmRange = par;\n    break;\n  }\n}\n\nif (ptmRange) {\n  ptmRange.insertText(\" (Kijin)\", \"End\");\n}\n\n// 2) Remove the trailing empty paragraph (lastRenderedPageBreak + tab, no text)\n// that follows \"Users bring in design questions\" at the very end of the body.\nconst allParas = body.paragraphs;\nallParas.load(\"items/text\");\nawait context.sync();\n\nconst lastPara = allParas.items[allParas.items.length - 1];\nif (lastPara && lastPara.text.trim() === \"\") {\n  lastPara.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Append \" (Kijin)\" to the end of the \"1:15  PTM\" line (Day 2 PTM tutorial item).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"1:15  PTM\")\nif ($found) {\n    $rng.Collapse(0)\n    $rng.InsertAfter(\" (Kijin)\")\n}\n\n# 2) Remove the trailing empty paragraph (lastRenderedPageBreak + tab, no visible text)\n# that follows \"Users bring in design questions\" at the very end of the document.\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count)\nif ($lastPara.Range.Text.Trim() -eq \"\") {\n    $lastPara.Range.Delete()\n}\n"}
